$wb = $excel.ActiveWorkbook

# Rows (in each localization sheet) that correspond to files that just had a
# handoff xliff re-generated for them.
$rows = @(7, 8, 9, 11, 12, 14)

# --- Overview sheet -------------------------------------------------------
# Column G = "Latest HO Xliff Generate Date"
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-09-04 00:24:33"
}

# --- zh-cn sheet ------------------------------------------------------------
# Column E = "Priority", Column H = "Latest Handoff Datetime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2016-09-04 00:24:29"
}

# --- de-de sheet ------------------------------------------------------------
# Column E = "Priority", Column H = "Latest Handoff Datetime"
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = "2016-09-04 00:24:33"
}
